$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45171 -> 45172, i.e. 2023-09-02 -> 2023-09-03) for every data row (2..307).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 307 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45172
